$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $value) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '70.118.90'
Set-TextValue "E2" '  +1.16%  '

# Row 3
Set-TextValue "D3" '3.502.62'
Set-TextValue "E3" '  +0.23%  '

# Row 4
Set-TextValue "E4" '  +0.02%  '

# Row 5
Set-TextValue "D5" '603.27'
Set-TextValue "E5" '  -0.39%  '

# Row 6
Set-TextValue "D6" '175.25'
Set-TextValue "E6" '  +3.76%  '

# Row 7
Set-TextValue "D7" '0.611'
Set-TextValue "E7" '  -0.87%  '

# Row 8
Set-TextValue "D8" '3.496.27'
Set-TextValue "E8" '  +0.25%  '

# Row 9
Set-TextValue "E9" '  -0.03%  '

# Row 10
Set-TextValue "E10" '  -0.11%  '

# Row 11
Set-TextValue "D11" '7.21'
Set-TextValue "E11" '  +8.52%  '

# Row 12
Set-TextValue "E12" '  +0.73%  '

# Row 13
Set-TextValue "D13" '46.12'
Set-TextValue "E13" '  -1.25%  '

# Row 14
Set-TextValue "D14" '0.0000274'
Set-TextValue "E14" '  -0.60%  '

# Row 15
Set-TextValue "D15" '4.058.02'
Set-TextValue "E15" '  +0.10%  '

# Row 16
Set-TextValue "D16" '8.28'
Set-TextValue "E16" '  +0.16%  '

# Row 17
Set-TextValue "D17" '609.95'
Set-TextValue "E17" '  +0.14%  '

# Row 18
Set-TextValue "D18" '3.500.48'
Set-TextValue "E18" '  +0.22%  '

# Row 19
Set-TextValue "D19" '70.193.33'
Set-TextValue "E19" '  +1.25%  '

# Row 20
Set-TextValue "E20" '  +1.04%  '

# Row 21
Set-TextValue "D21" '17.30'
Set-TextValue "E21" '  +0.92%  '

# Row 22
Set-TextValue "D22" '0.874'
Set-TextValue "E22" '  +0.00%  '

# Row 23
Set-TextValue "D23" '8.97'
Set-TextValue "E23" '  -9.89%  '

# Row 24
Set-TextValue "D24" '98.27'
Set-TextValue "E24" '  +3.13%  '

# Row 25
Set-TextValue "E25" '  -1.20%  '

# Row 26
Set-TextValue "E26" '  -3.30%  '

# Row 27
Set-TextValue "E27" '  +0.09%  '

# Row 28
Set-TextValue "E28" '  -1.13%  '

# Row 29
Set-TextValue "D29" '33.82'
Set-TextValue "E29" '  +2.36%  '

# Row 30
Set-TextValue "D30" '8.98'
Set-TextValue "E30" '  -2.08%  '

# Row 31
Set-TextValue "B31" 'Filecoin'
Set-TextValue "C31" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D31" '8.01'
Set-TextValue "E31" '  -4.37%  '

# Row 32
Set-TextValue "B32" 'Stacks'
Set-TextValue "C32" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D32" '2.95'
Set-TextValue "E32" '  -3.59%  '

# Row 33
Set-TextValue "B33" 'Mantle'
Set-TextValue "C33" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D33" '1.28'
Set-TextValue "E33" '  -3.59%  '

# Row 34
Set-TextValue "B34" 'Bittensor'
Set-TextValue "C34" 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue "D34" '631.87'
Set-TextValue "E34" '  +14.12%  '

# Row 35
Set-TextValue "D35" '6.81'
Set-TextValue "E35" '  -0.71%  '

# Row 36
Set-TextValue "E36" '  -1.24%  '

# Row 37
Set-TextValue "B37" 'Cosmos'
Set-TextValue "C37" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D37" '10.72'
Set-TextValue "E37" '  +0.14%  '

# Row 38
Set-TextValue "B38" 'dogwifhat'
Set-TextValue "C38" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D38" '3.53'
Set-TextValue "E38" '  +2.10%  '

# Row 39
Set-TextValue "D39" '0.0472'
Set-TextValue "E39" '  +6.26%  '

# Row 40
Set-TextValue "D40" '56.65'
Set-TextValue "E40" '  +0.11%  '

# Row 41
Set-TextValue "E41" '  -0.05%  '

# Row 42
Set-TextValue "D42" '0.141'
Set-TextValue "E42" '  +2.64%  '

# Row 43
Set-TextValue "D43" '3.359.99'
Set-TextValue "E43" '  +0.99%  '

# Row 44
Set-TextValue "D44" '0.0₃0730'
Set-TextValue "E44" '  +5.33%  '

# Row 45
Set-TextValue "E45" '  -4.86%  '

# Row 46
Set-TextValue "E46" '  -2.38%  '

# Row 47
Set-TextValue "E47" '  +0.86%  '

# Row 48
Set-TextValue "D48" '2.55'
Set-TextValue "E48" '  -1.41%  '

# Row 49
Set-TextValue "E49" '  +0.82%  '

# Row 50
Set-TextValue "D50" '132.64'
Set-TextValue "E50" '  -1.84%  '

# Row 51
Set-TextValue "E51" '  -0.01%  '
